# Scheduled runner update: refresh market-derived profit figures across
# all eight job sheets (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR) in the
# Cactuar_Profits workbook. Values below are the refreshed
# currentAveragePrice* / LevePrice* / LeveProfit* figures for the specific
# leves whose market data changed since the last run.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 181.2
$ws.Range("I33").Value = 215.53334
$ws.Range("J33").Value = 78.2
$ws.Range("K33").Value = 215.53334
$ws.Range("L33").Value = 78.2
$ws.Range("M33").Value = 13.46665999999999
$ws.Range("N33").Value = -536.2
$ws.Range("H40").Value = 33532.715
$ws.Range("I40").Value = 27932.25
$ws.Range("J40").Value = 41000
$ws.Range("K40").Value = 27932.25
$ws.Range("L40").Value = 41000
$ws.Range("M40").Value = -27757.25
$ws.Range("N40").Value = -41350
$ws.Range("H62").Value = 4362.3184
$ws.Range("J62").Value = 5666.3335
$ws.Range("L62").Value = 5666.3335
$ws.Range("N62").Value = -6914.3335
$ws.Range("H65").Value = 4362.3184
$ws.Range("J65").Value = 5666.3335
$ws.Range("L65").Value = 28331.6675
$ws.Range("N65").Value = -34571.6675
$ws.Range("H112").Value = 3069.0212
$ws.Range("J112").Value = 3086.848
$ws.Range("L112").Value = 9260.544
$ws.Range("N112").Value = -11476.544
$ws.Range("H132").Value = 144581.12
$ws.Range("I132").Value = 603477.0600000001
$ws.Range("K132").Value = 1810431.18
$ws.Range("M132").Value = -1807901.18
$ws.Range("H137").Value = 22225186
$ws.Range("I137").Value = 1757
$ws.Range("J137").Value = 41670690
$ws.Range("K137").Value = 5271
$ws.Range("L137").Value = 125012070
$ws.Range("M137").Value = -2721
$ws.Range("N137").Value = -125017170
$ws.Range("H138").Value = 5855.9473
$ws.Range("I138").Value = 2340.6155
$ws.Range("J138").Value = 6581.3335
$ws.Range("K138").Value = 7021.8465
$ws.Range("L138").Value = 19744.0005
$ws.Range("M138").Value = -1881.8465
$ws.Range("N138").Value = -30024.0005
$ws.Range("H140").Value = 68303.336
$ws.Range("I140").Value = 70000
$ws.Range("J140").Value = 68091.25
$ws.Range("K140").Value = 70000
$ws.Range("L140").Value = 68091.25
$ws.Range("M140").Value = -64820
$ws.Range("N140").Value = -78451.25

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 11637849
$ws.Range("I61").Value = 16137230
$ws.Range("K61").Value = 16137230
$ws.Range("M61").Value = -16137018
$ws.Range("H74").Value = 27782604
$ws.Range("I74").Value = 75001380
$ws.Range("K74").Value = 75001380
$ws.Range("M74").Value = -75000506
$ws.Range("H77").Value = 27782604
$ws.Range("I77").Value = 75001380
$ws.Range("K77").Value = 375006900
$ws.Range("M77").Value = -375002532
$ws.Range("H136").Value = 11637849
$ws.Range("I136").Value = 16137230
$ws.Range("K136").Value = 48411690
$ws.Range("M136").Value = -48409140
$ws.Range("H140").Value = 111999.664
$ws.Range("J140").Value = 111999.664
$ws.Range("L140").Value = 111999.664
$ws.Range("N140").Value = -122359.664

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H99").Value = 2606291
$ws.Range("I99").Value = 4168665.5
$ws.Range("K99").Value = 4168665.5
$ws.Range("M99").Value = -4167167.5
$ws.Range("H105").Value = 2141.0833
$ws.Range("I105").Value = 1829.8572
$ws.Range("K105").Value = 1829.8572
$ws.Range("M105").Value = -82.85719999999992
$ws.Range("H107").Value = 31251702
$ws.Range("I107").Value = 41668550
$ws.Range("K107").Value = 41668550
$ws.Range("M107").Value = -41666630
$ws.Range("H134").Value = 1361.625
$ws.Range("I134").Value = 1149.0834
$ws.Range("K134").Value = 3447.2502
$ws.Range("M134").Value = -912.2501999999999
$ws.Range("H140").Value = 121681.21
$ws.Range("I140").Value = 51000
$ws.Range("J140").Value = 140957.9
$ws.Range("K140").Value = 51000
$ws.Range("L140").Value = 140957.9
$ws.Range("M140").Value = -45820
$ws.Range("N140").Value = -151317.9

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 4850.5186
$ws.Range("J31").Value = 7086.7393
$ws.Range("L31").Value = 7086.7393
$ws.Range("N31").Value = -7676.7393
$ws.Range("H34").Value = 4850.5186
$ws.Range("J34").Value = 7086.7393
$ws.Range("L34").Value = 7086.7393
$ws.Range("N34").Value = -7490.7393
$ws.Range("H58").Value = 558126.9399999999
$ws.Range("I58").Value = 771570.5600000001
$ws.Range("K58").Value = 771570.5600000001
$ws.Range("M58").Value = -771367.5600000001
$ws.Range("H106").Value = 41999
$ws.Range("I106").Value = 0
$ws.Range("K106").Value = 0
$ws.Range("M106").ClearContents()
$ws.Range("H107").Value = 2256.5217
$ws.Range("I107").Value = 392.3
$ws.Range("K107").Value = 392.3
$ws.Range("M107").Value = 1527.7
$ws.Range("H132").Value = 1430.3334
$ws.Range("I132").Value = 734.25
$ws.Range("K132").Value = 2202.75
$ws.Range("M132").Value = 327.25
$ws.Range("H136").Value = 558126.9399999999
$ws.Range("I136").Value = 771570.5600000001
$ws.Range("K136").Value = 2314711.68
$ws.Range("M136").Value = -2312161.68
$ws.Range("H141").Value = 82281
$ws.Range("J141").Value = 87651.87
$ws.Range("L141").Value = 87651.87
$ws.Range("N141").Value = -98011.87

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 2009.8481
$ws.Range("J131").Value = 1982.6567
$ws.Range("L131").Value = 5947.9701
$ws.Range("N131").Value = -16027.9701

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H140").Value = 0
$ws.Range("J140").Value = 0
$ws.Range("L140").Value = 0
$ws.Range("N140").ClearContents()

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 7185
$ws.Range("J46").Value = 7822.143
$ws.Range("L46").Value = 7822.143
$ws.Range("N46").Value = -8198.143
$ws.Range("H61").Value = 2822.8462
$ws.Range("J61").Value = 2474.25
$ws.Range("L61").Value = 2474.25
$ws.Range("N61").Value = -2878.25
$ws.Range("H113").Value = 2822.8462
$ws.Range("J113").Value = 2474.25
$ws.Range("L113").Value = 2474.25
$ws.Range("N113").Value = -6814.25
$ws.Range("H136").Value = 5433.156
$ws.Range("I136").Value = 4987.086
$ws.Range("K136").Value = 14961.258
$ws.Range("M136").Value = -12411.258
$ws.Range("H139").Value = 85000
$ws.Range("J139").Value = 85000
$ws.Range("L139").Value = 85000
$ws.Range("N139").Value = -95280

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 7775
$ws.Range("I62").Value = 6480
$ws.Range("K62").Value = 6480
$ws.Range("M62").Value = -5856
$ws.Range("H65").Value = 7775
$ws.Range("I65").Value = 6480
$ws.Range("K65").Value = 32400
$ws.Range("M65").Value = -29280
$ws.Range("H100").Value = 871536.1
$ws.Range("I100").Value = 1177975
$ws.Range("J100").Value = 3292.5
$ws.Range("K100").Value = 2355950
$ws.Range("L100").Value = 6585
$ws.Range("M100").Value = -2355409
$ws.Range("N100").Value = -7667
$ws.Range("H136").Value = 8583.906000000001
$ws.Range("I136").Value = 2695.9333
$ws.Range("K136").Value = 8087.7999
$ws.Range("M136").Value = -5537.7999
